$d = $word.ActiveDocument

# The document ends with the "...Chris' explanation..." paragraph, followed
# by two pre-existing empty paragraphs, then the sectPr. The diff inserts a
# new "Caveat" block of 14 paragraphs right after those two empty
# paragraphs (i.e. immediately before the very last paragraph of the
# document). We anchor on the last paragraph and InsertXML a run of new
# <w:p> elements there: InsertXML folds the original trailing paragraph
# mark into the *last* <w:p> we supply, so we include one extra empty
# paragraph at the end of our fragment to stand in for it, which keeps the
# pre-existing final (now second-to-last) paragraph completely untouched.

$last = $d.Paragraphs.Last
$anchor = $d.Range($last.Range.Start, $last.Range.Start)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Caveat – Weird interaction</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Plotting bars with activations ends up showing that</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>– rACC is only activated in the R condition</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>– amygdalas are activated in both R and C</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>– fusiform is activated in both R and C</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>which is weird. my intuition is that this effect in the t-statistics is driven by the number of trials</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(R) = 2873</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(L) = 797</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(I) = 741</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(C) = 1240</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>but chris says that can't be the reason.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr></w:p>
'@

$anchor.InsertXML($xml)
